$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 378, shifting the existing row 378 (and everything
# below it, through the old row 409) down to become row 379 (through 410).
$ws.Rows.Item(378).Insert()

# Populate the newly inserted row 378 with the new weekly price record.
$ws.Cells.Item(378, 1).Value = 3
$ws.Cells.Item(378, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(378, 3).Value = "Coquimbo"
$ws.Cells.Item(378, 4).Value = 44769
$ws.Cells.Item(378, 4).NumberFormat = $ws.Cells.Item(379, 4).NumberFormat
$ws.Cells.Item(378, 5).Value = 5
$ws.Cells.Item(378, 6).Value = 100112031
$ws.Cells.Item(378, 7).Value = "Poroto verde"
$ws.Cells.Item(378, 8).Value = "Magnum"
$ws.Cells.Item(378, 9).Value = "Primera"
$ws.Cells.Item(378, 10).Value = 85
$ws.Cells.Item(378, 11).Value = 32000
$ws.Cells.Item(378, 12).Value = 33000
$ws.Cells.Item(378, 13).Value = 32471
$ws.Cells.Item(378, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(378, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(378, 16).Value = 1299
$ws.Cells.Item(378, 17).Value = 25
$ws.Cells.Item(378, 18).Value = "Hortaliza"
